$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old "intervention_type" marker / label values to their new values.
# A (column 1): colored-square glyph used as a visual marker
# B (column 2): the French color-name label describing the marker
$markerMap = @{
    "⬛" = "📘"
    "🟥" = "📕"
    "🟧" = "📙"
    "🟩" = "📗"
}
$labelMap = @{
    "noir" = "bleu"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valueA = $cellA.Value2
    if ($markerMap.ContainsKey($valueA)) {
        $cellA.Value2 = $markerMap[$valueA]
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valueB = $cellB.Value2
    if ($labelMap.ContainsKey($valueB)) {
        $cellB.Value2 = $labelMap[$valueB]
    }
}
